# Update "想去人数" (want-to-go count) figures in column F across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 149
$ws1.Range("F3").Value = 1795
$ws1.Range("F4").Value = 37
$ws1.Range("F6").Value = 668
$ws1.Range("F7").Value = 40
$ws1.Range("F13").Value = 163
$ws1.Range("F17").Value = 110
$ws1.Range("F18").Value = 5077
$ws1.Range("F19").Value = 56
$ws1.Range("F20").Value = 836
$ws1.Range("F22").Value = 2269
$ws1.Range("F24").Value = 29
$ws1.Range("F25").Value = 2112

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 83

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 149
$ws4.Range("F3").Value = 1795
$ws4.Range("F4").Value = 37
$ws4.Range("F6").Value = 668
$ws4.Range("F7").Value = 40
$ws4.Range("F13").Value = 163
$ws4.Range("F17").Value = 110
$ws4.Range("F18").Value = 5078
$ws4.Range("F19").Value = 83
$ws4.Range("F20").Value = 56
$ws4.Range("F22").Value = 836
$ws4.Range("F24").Value = 2269
$ws4.Range("F26").Value = 29
$ws4.Range("F27").Value = 2112
